$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (shifts existing rows 3-5 down to 4-6)
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the Bolivia match data
$ws.Cells.Item(3, 1).Value = 'Cn4YuSZh'
$ws.Cells.Item(3, 2).Value = '25/11/2024'
$ws.Cells.Item(3, 3).Value = '21:00'
$ws.Cells.Item(3, 4).Value = 'BOLIVIA - DIVISION PROFESIONAL'
$ws.Cells.Item(3, 5).Value = 'Oriente Petrolero'
$ws.Cells.Item(3, 6).Value = 'Santa Cruz'
$ws.Cells.Item(3, 7).Value = 1.48
$ws.Cells.Item(3, 8).Value = 4.1
$ws.Cells.Item(3, 9).Value = 7
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(3, 11).Value = 2.3
$ws.Cells.Item(3, 12).Value = 6.5
$ws.Cells.Item(3, 13).Value = 1.05
$ws.Cells.Item(3, 14).Value = 11
$ws.Cells.Item(3, 15).Value = 1.25
$ws.Cells.Item(3, 16).Value = 3.75
$ws.Cells.Item(3, 17).Value = 1.8
$ws.Cells.Item(3, 18).Value = 2
$ws.Cells.Item(3, 19).Value = 1.36
$ws.Cells.Item(3, 20).Value = 3
$ws.Cells.Item(3, 21).Value = 1.91
$ws.Cells.Item(3, 22).Value = 1.8
$ws.Cells.Item(3, 23).Value = 7
$ws.Cells.Item(3, 24).Value = 7
$ws.Cells.Item(3, 25).Value = 8.5
$ws.Cells.Item(3, 26).Value = 10
$ws.Cells.Item(3, 27).Value = 12
$ws.Cells.Item(3, 28).Value = 26
$ws.Cells.Item(3, 29).Value = 11
$ws.Cells.Item(3, 30).Value = 8
$ws.Cells.Item(3, 31).Value = 17
$ws.Cells.Item(3, 32).Value = 51
$ws.Cells.Item(3, 33).Value = 301
$ws.Cells.Item(3, 34).Value = 17
$ws.Cells.Item(3, 35).Value = 34
$ws.Cells.Item(3, 36).Value = 21
$ws.Cells.Item(3, 37).Value = 81
$ws.Cells.Item(3, 38).Value = 51
$ws.Cells.Item(3, 39).Value = 51
$ws.Cells.Item(3, 40).Value = 3.4
$ws.Cells.Item(3, 41).Value = 7.5
$ws.Cells.Item(3, 42).Value = 19
$ws.Cells.Item(3, 43).Value = 21
$ws.Cells.Item(3, 44).Value = 41
$ws.Cells.Item(3, 45).Value = 151
$ws.Cells.Item(3, 46).Value = 3
$ws.Cells.Item(3, 47).Value = 9
$ws.Cells.Item(3, 48).Value = 51
$ws.Cells.Item(3, 49).Value = 8
$ws.Cells.Item(3, 50).Value = 34
$ws.Cells.Item(3, 51).Value = 41
$ws.Cells.Item(3, 52).Value = 126
$ws.Cells.Item(3, 53).Value = 151
$ws.Cells.Item(3, 54).Value = 301
$ws.Cells.Item(3, 55).Value = 51
$ws.Cells.Item(3, 56).Value = 51

# Update odds that changed for the shifted rows (now rows 4, 5, 6)
# Row 4 (was row 3)
$ws.Cells.Item(4, 13).Value = 1.08
$ws.Cells.Item(4, 15).Value = 1.44
$ws.Cells.Item(4, 16).Value = 2.63
# Row 5 (was row 4)
$ws.Cells.Item(5, 13).Value = 1.1
$ws.Cells.Item(5, 15).Value = 1.44
$ws.Cells.Item(5, 16).Value = 2.63
$ws.Cells.Item(5, 17).Value = 2.5
$ws.Cells.Item(5, 18).Value = 1.5
# Row 6 (was row 5)
$ws.Cells.Item(6, 7).Value = 2.3
$ws.Cells.Item(6, 13).Value = 1.06
$ws.Cells.Item(6, 14).Value = 10
$ws.Cells.Item(6, 17).Value = 2.15
$ws.Cells.Item(6, 18).Value = 1.67
$ws.Cells.Item(6, 40).Value = 4.33
$ws.Cells.Item(6, 41).Value = 13
